# Updated cryptos list with refreshed prices and volume percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text is a plain decimal number need the "@" (Text)
# number format first, otherwise Excel auto-converts the typed text into a
# real number and silently drops meaningful trailing zeros (e.g. "20.20" -> 20.2).

$ws.Range("D2").Value = "69.276.14"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "2.752.21"
$ws.Range("E3").Value = "  +3.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.22"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.93"
$ws.Range("E6").Value = "  +5.19%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.548"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").Value = "2.752.84"
$ws.Range("E9").Value = "  +3.95%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  +3.78%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.34"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("D15").Value = "3.253.05"
$ws.Range("E15").Value = "  +3.91%  "
$ws.Range("D17").Value = "69.158.30"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").Value = "2.763.77"
$ws.Range("E18").Value = "  +3.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.98"
$ws.Range("E19").Value = "  +4.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.75"
$ws.Range("E20").Value = "  +5.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.06"
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.99"
$ws.Range("E23").Value = "  +3.53%  "
$ws.Range("E24").Value = "  +3.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.33"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.08"
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("D28").Value = "2.877.09"
$ws.Range("E28").Value = "  +3.33%  "
$ws.Range("E29").Value = "  +2.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "604.12"
$ws.Range("E30").Value = "  +8.10%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.36"
$ws.Range("E32").Value = "  +4.15%  "
$ws.Range("E33").Value = "  +3.71%  "
$ws.Range("E34").Value = "  +5.85%  "
$ws.Range("E35").Value = "  +3.71%  "
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.20"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "163.32"
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.74"
$ws.Range("E43").Value = "  +4.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.06"
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").Value = "0.0₆0320"
$ws.Range("E45").Value = "  -3.62%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.70"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.99"
$ws.Range("E48").Value = "  +5.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.81"
$ws.Range("E49").Value = "  +6.84%  "
$ws.Range("E50").Value = "  +8.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.26"
$ws.Range("E51").Value = "  +0.41%  "
